$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.75

# Row 4
$ws.Range("M4").Value = 1.07
$ws.Range("O4").Value = 1.36

# Row 5
$ws.Range("M5").Value = 1.05
$ws.Range("O5").Value = 1.29
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93

# Row 8
$ws.Range("G8").Value = 1.95
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 2.6
$ws.Range("Z8").Value = 17
$ws.Range("AH8").Value = 13
$ws.Range("AI8").Value = 21
$ws.Range("AN8").Value = 4
$ws.Range("AZ8").Value = 67
